# "Finished Week 13 logging" - update the Road ("R") row of target depth
# data on both the OFF and DEF sheets with the latest cumulative totals.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 549
$wsOff.Range("C3").Value = 399
$wsOff.Range("D3").Value = 155
$wsOff.Range("E3").Value = 75
$wsOff.Range("F3").Value = 6

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 579
$wsDef.Range("C3").Value = 420
$wsDef.Range("D3").Value = 100
$wsDef.Range("E3").Value = 54
